$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-15 Wednesday" "2025-01-16 Thursday"

Replace-Text "685×2=" "486×5="
Replace-Text "481×4=" "837×7="
Replace-Text "951×8=" "486×9="
Replace-Text "608×3=" "750×4="
Replace-Text "329×5=" "981×8="

Replace-Text "973×2=" "478×3="
Replace-Text "150×9=" "609×7="
Replace-Text "663×8=" "702×2="
Replace-Text "827×8=" "577×5="
Replace-Text "255×3=" "719×8="

Replace-Text "264×3=" "612×3="
Replace-Text "641×6=" "137×9="
Replace-Text "891×3=" "141×9="
Replace-Text "169×7=" "371×2="
Replace-Text "348×4=" "308×9="

Replace-Text "377×5=" "714×8="
Replace-Text "840×5=" "971×6="
Replace-Text "956×9=" "256×5="
Replace-Text "994×2=" "640×5="
Replace-Text "712×7=" "206×4="

Replace-Text "671×6=" "416×2="
Replace-Text "276×5=" "720×5="
Replace-Text "951×2=" "399×4="
Replace-Text "150×4=" "830×6="
Replace-Text "473×2=" "730×3="
